$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    "B2=18.08625015614324", "C2=9.73871745288309", "D2=7.997305569775516", "E2=12.98521812365697", "F2=36.9053411195708", "J2=10.22553441222674", "L2=10.92129140794003", "M2=16.81086850035017", "O2=28.45777371119366",
    "B3=17.64250469793256", "C3=9.515025990554125", "D3=7.999093464023402", "E3=13.02009814429155", "F3=37.04063133546499", "J3=10.24710118614294", "L3=10.91766058614509", "M3=16.70607724895474", "O3=28.57547145604947",
    "B4=17.36624223006073", "C4=9.374347221146873", "D4=8.000850131172871", "E4=13.04269451355855", "F4=37.13270598905552", "J4=10.26102417996157", "L4=10.91645423496688", "M4=16.64321381675922", "O4=28.65438626626341",
    "B5=17.25286273582088", "C5=9.316240138901346", "D5=8.001732233848356", "E5=13.05220020832665", "F5=37.17248566416499", "J5=10.26686966170448", "L5=10.91622105425597", "M5=16.61798642928891", "O5=28.68821249252099",
    "B6=17.23399251259568", "C6=9.306546200830589", "D6=8.001888763253007", "E6=13.05379661286609", "F6=37.17922727479164", "J6=10.26785068826055", "L6=10.91619797577643", "M6=16.61382151245772", "O6=28.69392992849943",
    "B7=17.36471618393417", "C7=9.373566644748927", "D7=8.000861353655706", "E7=13.04282150509225", "F7=37.13323333530997", "J7=10.26110231806805", "L7=10.91645004237555", "M7=16.64287198890554", "O7=28.65483571090578",
    "B8=17.93412649843779", "C8=9.662314023781489", "D8=7.997785624862069", "E8=12.99700034779674", "F8=36.95011639964408", "J8=10.2328296327119", "L8=10.91982790202491", "M8=16.77444057034848", "O8=28.4969738072995",
    "B9=19.01415583286876", "C9=10.19972357362075", "D9=7.996957290825884", "E9=12.91647073980007", "F9=36.66274751145447", "J9=10.18276543148545", "L9=10.93451580147221", "M9=17.04336286279484", "O9=28.24032636897383",
    "B10=19.77736692925929", "C10=10.57403431950749", "D10=7.999487087165324", "E10=12.86293971281622", "F10=36.49569197463768", "J10=10.14922801773621", "L10=10.95014996457146", "M10=17.24647410129744", "O10=28.08425919176366",
    "B11=20.11649095246739", "C11=10.73930793975749", "D11=8.001311853207616", "E11=12.83979961382783", "F11=36.42933460500375", "J11=10.13466829144124", "L11=10.95829715904492", "M11=17.33983062985391", "O11=28.02036374992895",
    "B12=20.24363613226985", "C12=10.80113195028817", "D12=8.002099062000861", "E12=12.83121044968214", "F12=36.40559796223498", "J12=10.12925452808307", "L12=10.96152957711771", "M12=17.37529927681542", "O12=27.99719315660433",
    "B13=20.21631150781793", "C13=10.78785154747384", "D13=8.001925255031688", "E13=12.83305257517602", "F13=36.41064810084909", "J13=10.13041605283715", "L13=10.96082689513922", "M13=17.36765562918291", "O13=28.00213769526479",
    "B14=20.12697741001235", "C14=10.74440972975689", "D14=8.001374693285971", "E14=12.8390895052064", "F14=36.42735385899771", "J14=10.13422090270696", "L14=10.95856015026634", "M14=17.3427464445441", "O14=28.01843692046555",
    "B15=20.07208859522193", "C15=10.71769999609094", "D15=8.001049966327995", "E15=12.84280987150908", "F15=36.43776798924588", "J15=10.13656445130448", "L15=10.95719083158065", "M15=17.32750338783007", "O15=28.02855430434786",
    "B16=19.75503213112887", "C16=10.56312895429647", "D16=7.99938133642392", "E16=12.86447628869599", "F16=36.5002229972138", "J16=10.15019350632226", "L16=10.94963822069049", "M16=17.2403905179219", "O16=28.08857814715375",
    "B17=19.55837856966173", "C17=10.46699229366429", "D17=7.998529807932293", "E17=12.87807770121591", "F17=36.54100951890884", "J17=10.15873256415321", "L17=10.94526894542923", "M17=17.18718026079899", "O17=28.12722241004072",
    "B18=19.44451578194957", "C18=10.41122849539194", "D18=7.99810351023703", "E18=12.88601494205505", "F18=36.56537561050362", "J18=10.16370960132922", "L18=10.94285338564984", "M18=17.156666806557", "O18=28.15011763520762",
    "B19=19.40583812493957", "C19=10.39226870382466", "D19=7.997970096753646", "E19=12.8887219685839", "F19=36.57378111055935", "J19=10.16540602172612", "L19=10.94205231739913", "M19=17.14635187447151", "O19=28.15798417662024",
    "B20=19.57939145150643", "C20=10.47727501363752", "D20=7.9986138897833", "E20=12.87661800720432", "F20=36.53657383547843", "J20=10.157816781338", "L20=10.94572397943479", "M20=17.1928352506028", "O20=28.12303948903075",
    "B21=20.15325243070728", "C21=10.75719063129957", "D21=8.001533801250527", "E21=12.83731160972661", "F21=36.42240916565144", "J21=10.13310062436605", "L21=10.95922196487213", "M21=17.35005988471357", "O21=28.01362158553802",
    "B22=20.52083321655463", "C22=10.93567578680518", "D22=8.004002570609982", "E22=12.81263350792647", "F22=36.3559088135361", "J22=10.11752802962513", "L22=10.96890118889087", "M22=17.4534845942982", "O22=27.9480883378817",
    "B23=20.32536723465662", "C23=10.84083566395278", "D23=8.00263390070678", "E23=12.82571241176063", "F23=36.39065718581939", "J23=10.12578642785771", "L23=10.96365728820521", "M23=17.39823077528673", "O23=27.98251634024198",
    "B24=19.56989402324673", "C24=10.47262773003472", "D24=7.998575679299334", "E24=12.8772775684758", "F24=36.53857635077588", "J24=10.1582305954555", "L24=10.94551795818425", "M24=17.19027838563054", "O24=28.12492847558235",
    "B25=18.72675575454393", "C25=10.0577444681157", "D25=7.996627718331036", "E25=12.93726306837049", "F25=36.73277286391109", "J25=10.19573686370757", "L25=10.92968630897144", "M25=16.9695577510326", "O25=28.3040668631896"
)

foreach ($entry in $updates) {
    $parts = $entry -split "="
    $ref = $parts[0]
    $val = [double]$parts[1]
    $ws.Range($ref).Value = $val
}
